$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 62, shifting existing rows 62-73 down to 63-74
$ws.Rows.Item(62).Insert()

# Populate the newly inserted row 62 with the new record
$ws.Cells.Item(62, 1).Value = 11
$ws.Cells.Item(62, 2).Value = "Vega Monumental Concepción"
$ws.Cells.Item(62, 3).Value = "Bíobío"
$ws.Cells.Item(62, 4).Value = 44474
$ws.Cells.Item(62, 5).Value = 8
$ws.Cells.Item(62, 6).Value = 100112043
$ws.Cells.Item(62, 7).Value = "Pepino ensalada"
$ws.Cells.Item(62, 8).Value = "Sin especificar"
$ws.Cells.Item(62, 9).Value = "Primera"
$ws.Cells.Item(62, 10).Value = 100
$ws.Cells.Item(62, 11).Value = 19000
$ws.Cells.Item(62, 12).Value = 20000
$ws.Cells.Item(62, 13).Value = 19500
$ws.Cells.Item(62, 14).Value = "`$/caja 60 unidades"
$ws.Cells.Item(62, 15).Value = "Región de Arica y Parinacota"
$ws.Cells.Item(62, 16).Value = 325
$ws.Cells.Item(62, 17).Value = 60
$ws.Cells.Item(62, 18).Value = "Hortaliza"
